# Add the new "2022-Q4" data: a new sheet with fund-holding detail, plus a new
# summary row at the top of the "总计" (total) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row 2 for 2022-Q4, pushing the
#    existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Grab the header's formatting (bold/centered/bordered) before the insert
# shifts things around - this is the same style used by column A below it.
$summary.Range("B1").Copy()

$summary.Range("A2").EntireRow.Insert()

# Re-apply the index-column style to A2:A7 (the insert only half-preserved
# it), then fill in the values.
$summary.Range("A2:A7").PasteSpecial(-4122)
# Data columns (B:D) should have no special style, same as the other rows.
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.79

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# 2. New "2022-Q4" sheet with the fund holdings for the quarter.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Add()
$q4.Name = "2022-Q4"

# Copy header formatting (bold / centered / bordered) from an existing sheet
# so the new sheet matches the rest of the workbook.
$template = $wb.Worksheets.Item("2022-Q3")
$template.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$q4.Range("A2:A6").PasteSpecial(-4122)

# The fund code / size / position columns hold numeric-looking values that
# must stay text (e.g. leading zeros in fund codes), same as in the other
# quarter sheets.
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "002418"
$q4.Range("C2").Value = "汇添富优选回报灵活配置混合C"
$q4.Range("D2").Value = "9.25"
$q4.Range("E2").Value = "94.49"
$q4.Range("F2").Value = "4.04"
$q4.Range("G2").Value = "0.3737"
$q4.Range("H2").Value = 10

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "470021"
$q4.Range("C3").Value = "汇添富优选回报灵活配置混合A"
$q4.Range("D3").Value = "5.37"
$q4.Range("E3").Value = "94.49"
$q4.Range("F3").Value = "4.04"
$q4.Range("G3").Value = "0.2169"
$q4.Range("H3").Value = 10

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "015696"
$q4.Range("C4").Value = "农银绿色能源混合"
$q4.Range("D4").Value = "3.40"
$q4.Range("E4").Value = "82.40"
$q4.Range("F4").Value = "3.48"
$q4.Range("G4").Value = "0.1183"
$q4.Range("H4").Value = 6

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "013250"
$q4.Range("C5").Value = "红土创新智能制造混合"
$q4.Range("D5").Value = "1.19"
$q4.Range("E5").Value = "90.41"
$q4.Range("F5").Value = "3.84"
$q4.Range("G5").Value = "0.0457"
$q4.Range("H5").Value = 8

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "004044"
$q4.Range("C6").Value = "金鹰转型动力灵活配置混合"
$q4.Range("D6").Value = "0.65"
$q4.Range("E6").Value = "90.51"
$q4.Range("F6").Value = "5.23"
$q4.Range("G6").Value = "0.0340"
$q4.Range("H6").Value = 7

# Move the new sheet right after "总计" so the tab order is:
# 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2021-Q4, 2021-Q3, 2021-Q2
$q4.Move($null, $summary)

Write-Host "Done adding 2022-Q4 sheet and summary row"
